$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatted "Reported Date" cell style ahead of writing new dates,
# so the new cells reuse the existing date-format style instead of Excel
# minting a brand-new numFmt/style entry.
$ws.Range("H3").Copy()
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null

# New row 5 — second query entry
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "1.5.0"
$ws.Range("C5").Value = "Pbook"
$ws.Range("D5").Value = "which fields need to be validate?"
$ws.Range("G5").Value = "Rosalin"
$ws.Range("H5").Value = [DateTime]"2020-06-25"
$ws.Range("J5").Value = "Ashok"
$ws.Range("K5").Value = "open"

# New row 7 — third query entry
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "1.5.0"
$ws.Range("C7").Value = "Pbook"
$ws.Range("D7").Value = "how many records should be loaded for once?"
$ws.Range("F7").Value = " "
$ws.Range("G7").Value = "Rosalin"
$ws.Range("H7").Value = [DateTime]"2020-06-25"
$ws.Range("J7").Value = "Ashok"
$ws.Range("K7").Value = "open"

$ws.Range("E15").Select()
